{"js": "// Office.js (Word JavaScript API) edit script.\n// Updates the date line and the 26 two-digit multiplication problems in\n// the document's table cells, replacing each old \"NN\u00d7NN=\" (or date)\n// string with its new value, per the commit's diff.\n\nconst replacements = [\n  [\"2025-12-19 Friday\", \"2025-12-20 Saturday\"],\n  [\"33\u00d737=\", \"78\u00d751=\"],\n  [\"54\u00d725=\", \"85\u00d767=\"],\n  [\"14\u00d749=\", \"36\u00d738=\"],\n  [\"35\u00d718=\", \"40\u00d742=\"],\n  [\"44\u00d798=\", \"14\u00d780=\"],\n  [\"74\u00d734=\", \"16\u00d733=\"],\n  [\"82\u00d761=\", \"20\u00d798=\"],\n  [\"13\u00d798=\", \"17\u00d777=\"],\n  [\"96\u00d772=\", \"81\u00d772=\"],\n  [\"46\u00d786=\", \"54\u00d715=\"],\n  [\"75\u00d734=\", \"30\u00d757=\"],\n  [\"87\u00d720=\", \"96\u00d778=\"],\n  [\"14\u00d715=\", \"55\u00d778=\"],\n  [\"49\u00d727=\", \"97\u00d727=\"],\n  [\"28\u00d740=\", \"49\u00d775=\"],\n  [\"90\u00d753=\", \"17\u00d775=\"],\n  [\"67\u00d711=\", \"95\u00d728=\"],\n  [\"55\u00d764=\", \"75\u00d713=\"],\n  [\"63\u00d793=\", \"30\u00d743=\"],\n  [\"48\u00d756=\", \"57\u00d732=\"],\n  [\"92\u00d724=\", \"50\u00d790=\"],\n  [\"60\u00d758=\", \"72\u00d726=\"],\n  [\"60\u00d729=\", \"73\u00d737=\"],\n  [\"98\u00d736=\", \"46\u00d782=\"],\n  [\"79\u00d757=\", \"96\u00d785=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Updates the date line and the 26 two-digit multiplication problems in\n# the document's table cells, replacing each old \"NN\u00d7NN=\" (or date)\n# string with its new value, per the commit's diff.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-19 Friday\", \"2025-12-20 Saturday\"),\n    @(\"33\u00d737=\", \"78\u00d751=\"),\n    @(\"54\u00d725=\", \"85\u00d767=\"),\n    @(\"14\u00d749=\", \"36\u00d738=\"),\n    @(\"35\u00d718=\", \"40\u00d742=\"),\n    @(\"44\u00d798=\", \"14\u00d780=\"),\n    @(\"74\u00d734=\", \"16\u00d733=\"),\n    @(\"82\u00d761=\", \"20\u00d798=\"),\n    @(\"13\u00d798=\", \"17\u00d777=\"),\n    @(\"96\u00d772=\", \"81\u00d772=\"),\n    @(\"46\u00d786=\", \"54\u00d715=\"),\n    @(\"75\u00d734=\", \"30\u00d757=\"),\n    @(\"87\u00d720=\", \"96\u00d778=\"),\n    @(\"14\u00d715=\", \"55\u00d778=\"),\n    @(\"49\u00d727=\", \"97\u00d727=\"),\n    @(\"28\u00d740=\", \"49\u00d775=\"),\n    @(\"90\u00d753=\", \"17\u00d775=\"),\n    @(\"67\u00d711=\", \"95\u00d728=\"),\n    @(\"55\u00d764=\", \"75\u00d713=\"),\n    @(\"63\u00d793=\", \"30\u00d743=\"),\n    @(\"48\u00d756=\", \"57\u00d732=\"),\n    @(\"92\u00d724=\", \"50\u00d790=\"),\n    @(\"60\u00d758=\", \"72\u00d726=\"),\n    @(\"60\u00d729=\", \"73\u00d737=\"),\n    @(\"98\u00d736=\", \"46\u00d782=\"),\n    @(\"79\u00d757=\", \"96\u00d785=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
